$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G (MAXHP), shifting MAXMP and everything after it one column right.
$ws.Columns("G:G").Insert()
$ws.Range("G1").ColumnWidth = 13.29

# New header for the inserted column.
$ws.Range("G1").Value = "MAXHP"

# NPC HP fix: MAXHP (new col G) and MAXMP (shifted col H) should equal the SalePrice (col F) value,
# and DEF_FIRE (shifted col R) goes from 5 to 10.
$ws.Range("G2").Value = $ws.Range("F2").Value2
$ws.Range("G3").Value = $ws.Range("F3").Value2
$ws.Range("G4").Value = $ws.Range("F4").Value2
$ws.Range("G5").Value = $ws.Range("F5").Value2
$ws.Range("G6").Value = $ws.Range("F6").Value2

$ws.Range("H2").Value = $ws.Range("F2").Value2
$ws.Range("H3").Value = $ws.Range("F3").Value2
$ws.Range("H4").Value = $ws.Range("F4").Value2
$ws.Range("H5").Value = $ws.Range("F5").Value2
$ws.Range("H6").Value = $ws.Range("F6").Value2

$ws.Range("R2").Value = 10
$ws.Range("R3").Value = 10
$ws.Range("R4").Value = 10
$ws.Range("R5").Value = 10
$ws.Range("R6").Value = 10

# Restore the view: scroll back to the top-left and select I8 (matches the saved workbook state).
$ws.Activate()
$ws.Range("I8").Select()
